# Refresh the cryptos table (coinranking.com price/%-change scrape)
# with the values captured by this run's GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'34.402.60"
$ws.Range("E2").Value = '  +0.85%  '

# Row 3
$ws.Range("D3").Value = "'1.795.79"
$ws.Range("E3").Value = '  +0.51%  '

# Row 4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").Value = "'226.71"
$ws.Range("E5").Value = '  +0.04%  '

# Row 6
$ws.Range("D6").Value = "'0.555"
$ws.Range("E6").Value = '  +1.48%  '

# Row 7
$ws.Range("E7").Value = '  +0.03%  '

# Row 8
$ws.Range("D8").Value = "'32.43"
$ws.Range("E8").Value = '  +1.50%  '

# Row 9
$ws.Range("E9").Value = '  +1.33%  '

# Row 10
$ws.Range("D10").Value = "'0.0693"
$ws.Range("E10").Value = '  +0.30%  '

# Row 11
$ws.Range("D11").Value = "'0.0950"
$ws.Range("E11").Value = '  +0.65%  '

# Row 12
$ws.Range("D12").Value = "'2.055.03"
$ws.Range("E12").Value = '  +0.54%  '

# Row 13
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").Value = "'11.10"
$ws.Range("E13").Value = '  -1.22%  '

# Row 14
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = "'1.801.57"
$ws.Range("E14").Value = '  +0.76%  '

# Row 15
$ws.Range("D15").Value = "'0.630"
$ws.Range("E15").Value = '  +1.70%  '

# Row 16
$ws.Range("D16").Value = "'34.366.64"
$ws.Range("E16").Value = '  +0.94%  '

# Row 17
$ws.Range("E17").Value = '  +0.64%  '

# Row 18
$ws.Range("D18").Value = "'68.40"
$ws.Range("E18").Value = '  +0.60%  '

# Row 19
$ws.Range("D19").Value = "'0.0₃0803"
$ws.Range("E19").Value = '  +3.17%  '

# Row 20
$ws.Range("D20").Value = "'246.60"
$ws.Range("E20").Value = '  +0.34%  '

# Row 21
$ws.Range("D21").Value = "'11.02"
$ws.Range("E21").Value = '  +1.79%  '

# Row 22
$ws.Range("E22").Value = '  +0.03%  '

# Row 23
$ws.Range("D23").Value = "'4.16"
$ws.Range("E23").Value = '  +1.74%  '

# Row 24
$ws.Range("D24").Value = "'2.07"
$ws.Range("E24").Value = '  +1.25%  '

# Row 25
$ws.Range("D25").Value = "'163.03"
$ws.Range("E25").Value = '  +0.86%  '

# Row 26
$ws.Range("E26").Value = '  +0.62%  '

# Row 27
$ws.Range("D27").Value = "'16.40"
$ws.Range("E27").Value = '  +0.50%  '

# Row 28
$ws.Range("E28").Value = '  +2.07%  '

# Row 29
$ws.Range("E29").Value = '  +0.04%  '

# Row 30
$ws.Range("B30").Value = 'Hedera'
$ws.Range("C30").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D30").Value = "'0.0522"
$ws.Range("E30").Value = '  +0.85%  '

# Row 31
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = "'1.23"
$ws.Range("E31").Value = '  -0.18%  '

# Row 32
$ws.Range("D32").Value = "'3.90"
$ws.Range("E32").Value = '  +8.06%  '

# Row 34
$ws.Range("E34").Value = '  +1.24%  '

# Row 35
$ws.Range("D35").Value = "'1.441.65"
$ws.Range("E35").Value = '  -0.64%  '

# Row 36
$ws.Range("E36").Value = '  +9.34%  '

# Row 37
$ws.Range("D37").Value = "'0.665"
$ws.Range("E37").Value = '  +3.04%  '

# Row 38
$ws.Range("E38").Value = '  +1.95%  '

# Row 39
$ws.Range("E39").Value = '  -0.95%  '

# Row 40
$ws.Range("D40").Value = "'83.63"
$ws.Range("E40").Value = '  +4.36%  '

# Row 41
$ws.Range("E41").Value = '  +1.34%  '

# Row 42
$ws.Range("E42").Value = '  +1.61%  '

# Row 43
$ws.Range("E43").Value = '  +2.84%  '

# Row 44
$ws.Range("E44").Value = '  +2.39%  '

# Row 45
$ws.Range("D45").Value = "'0.0524"
$ws.Range("E45").Value = '  +3.18%  '

# Row 46
$ws.Range("D46").Value = "'6.09"
$ws.Range("E46").Value = '  +0.90%  '

# Row 47
$ws.Range("E47").Value = '  +0.05%  '

# Row 48
$ws.Range("D48").Value = "'1.951.93"
$ws.Range("E48").Value = '  +0.34%  '

# Row 49
$ws.Range("D49").Value = "'105.70"
$ws.Range("E49").Value = '  -1.75%  '

# Row 50
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = "'0.0₆0130"
$ws.Range("E50").Value = '  -4.87%  '

# Row 51
$ws.Range("B51").Value = 'PaxDollar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D51").Value = "'1.00"
$ws.Range("E51").Value = '  +0.02%  '
